$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1981.875
$ws.Range("I4").Value = 1576.1666
$ws.Range("J4").Value = 3199
$ws.Range("K4").Value = 1576.1666
$ws.Range("L4").Value = 3199
$ws.Range("M4").Value = -1462.1666

$ws.Range("H6").Value = 174.71428
$ws.Range("I6").Value = 200.5
$ws.Range("J6").Value = 20
$ws.Range("K6").Value = 601.5
$ws.Range("L6").Value = 60
$ws.Range("M6").Value = -489.5
$ws.Range("N6").Value = -284

$ws.Range("H15").Value = 964727.6
$ws.Range("I15").Value = 964727.6
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 2894182.8
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -2894013.8

$ws.Range("H33").Value = 3565.3333
$ws.Range("I33").Value = 1856.75
$ws.Range("J33").Value = 6982.5
$ws.Range("K33").Value = 1856.75
$ws.Range("L33").Value = 6982.5
$ws.Range("M33").Value = -1627.75

$ws.Range("H70").Value = 2179.6667
$ws.Range("I70").Value = 989
$ws.Range("J70").Value = 2775
$ws.Range("K70").Value = 2967
$ws.Range("L70").Value = 8325
$ws.Range("M70").Value = -2697
$ws.Range("N70").Value = -8865

$ws.Range("H73").Value = 2179.6667
$ws.Range("I73").Value = 989
$ws.Range("J73").Value = 2775
$ws.Range("K73").Value = 2967
$ws.Range("L73").Value = 8325
$ws.Range("M73").Value = -2031
$ws.Range("N73").Value = -10197

$ws.Range("H86").Value = 10158.25
$ws.Range("I86").Value = 9413.200000000001
$ws.Range("J86").Value = 11400
$ws.Range("K86").Value = 9413.200000000001
$ws.Range("L86").Value = 11400
$ws.Range("M86").Value = -8290.200000000001

$ws.Range("H89").Value = 10158.25
$ws.Range("I89").Value = 9413.200000000001
$ws.Range("J89").Value = 11400
$ws.Range("K89").Value = 47066
$ws.Range("L89").Value = 57000
$ws.Range("M89").Value = -41450

$ws.Range("H111").Value = 74614.64
$ws.Range("I111").Value = 1470.7778
$ws.Range("J111").Value = 206273.6
$ws.Range("K111").Value = 4412.3334
$ws.Range("L111").Value = 618820.8
$ws.Range("M111").Value = -1345.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 545.75
$ws.Range("I4").Value = 395.5
$ws.Range("J4").Value = 696
$ws.Range("K4").Value = 395.5
$ws.Range("L4").Value = 696
$ws.Range("M4").Value = -279.5

$ws.Range("H8").Value = 1750
$ws.Range("I8").Value = 1750
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1750
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1606

$ws.Range("H32").Value = 7180.2104
$ws.Range("I32").Value = 2938.6
$ws.Range("J32").Value = 56665.668
$ws.Range("K32").Value = 2938.6
$ws.Range("L32").Value = 56665.668
$ws.Range("M32").Value = -2651.6

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H122").Value = 2569.4783
$ws.Range("I122").Value = 1867.2667
$ws.Range("J122").Value = 3886.125
$ws.Range("K122").Value = 5601.800099999999
$ws.Range("L122").Value = 11658.375
$ws.Range("M122").Value = -3151.800099999999

$ws.Range("H133").Value = 110551
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 110551
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 110551
$ws.Range("N133").Value = -115611

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 9500
$ws.Range("I10").Value = 9500
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 9500
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -9360

$ws.Range("H86").Value = 2767.5386
$ws.Range("I86").Value = 1007.6
$ws.Range("J86").Value = 5167.4546
$ws.Range("K86").Value = 1007.6
$ws.Range("L86").Value = 5167.4546
$ws.Range("M86").Value = 115.4
$ws.Range("N86").Value = -7413.4546

$ws.Range("H89").Value = 2767.5386
$ws.Range("I89").Value = 1007.6
$ws.Range("J89").Value = 5167.4546
$ws.Range("K89").Value = 5038
$ws.Range("L89").Value = 25837.273
$ws.Range("M89").Value = 578
$ws.Range("N89").Value = -37069.273

$ws.Range("H94").Value = 721.4375
$ws.Range("I94").Value = 643.6799999999999
$ws.Range("J94").Value = 999.1429000000001
$ws.Range("K94").Value = 643.6799999999999
$ws.Range("L94").Value = 999.1429000000001
$ws.Range("M94").Value = -192.6799999999999
$ws.Range("N94").Value = -1901.1429

$ws.Range("H134").Value = 5764
$ws.Range("I134").Value = 3147.125
$ws.Range("J134").Value = 10997.75
$ws.Range("K134").Value = 9441.375
$ws.Range("L134").Value = 32993.25
$ws.Range("M134").Value = -6906.375
$ws.Range("N134").Value = -38063.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 573.6923
$ws.Range("I7").Value = 604.1111
$ws.Range("J7").Value = 505.25
$ws.Range("K7").Value = 604.1111
$ws.Range("L7").Value = 505.25
$ws.Range("M7").Value = -491.1111

$ws.Range("H13").Value = 42000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 42000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 42000
$ws.Range("N13").Value = -42278

$ws.Range("H16").Value = 1077.75
$ws.Range("I16").Value = 937
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 937
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -650

$ws.Range("H31").Value = 7472.3696
$ws.Range("I31").Value = 2682.318
$ws.Range("J31").Value = 11863.25
$ws.Range("K31").Value = 2682.318
$ws.Range("L31").Value = 11863.25
$ws.Range("M31").Value = -2387.318
$ws.Range("N31").Value = -12453.25

$ws.Range("H34").Value = 7472.3696
$ws.Range("I34").Value = 2682.318
$ws.Range("J34").Value = 11863.25
$ws.Range("K34").Value = 2682.318
$ws.Range("L34").Value = 11863.25
$ws.Range("M34").Value = -2480.318
$ws.Range("N34").Value = -12267.25

$ws.Range("H107").Value = 2043.5769
$ws.Range("I107").Value = 1742.3334
$ws.Range("J107").Value = 3308.8
$ws.Range("K107").Value = 1742.3334
$ws.Range("L107").Value = 3308.8
$ws.Range("M107").Value = 177.6666

$ws.Range("H113").Value = 1077.75
$ws.Range("I113").Value = 937
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 937
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1233

$ws.Range("H123").Value = 49999
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 49999
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 49999
$ws.Range("N123").Value = -59799

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1500
$ws.Range("I51").Value = 1500
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 4500
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -4040
$ws.Range("N51").ClearContents()

$ws.Range("H131").Value = 5082.5625
$ws.Range("I131").Value = 6365
$ws.Range("J131").Value = 4899.357
$ws.Range("K131").Value = 19095
$ws.Range("L131").Value = 14698.071
$ws.Range("M131").Value = -14055
$ws.Range("N131").Value = -24778.071

$ws.Range("H137").Value = 2430.1538
$ws.Range("I137").Value = 1600
$ws.Range("J137").Value = 5197.3335
$ws.Range("K137").Value = 4800
$ws.Range("L137").Value = 15592.0005
$ws.Range("M137").Value = 300
$ws.Range("N137").Value = -25792.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3028.55
$ws.Range("I80").Value = 2889.1667
$ws.Range("J80").Value = 3237.625
$ws.Range("K80").Value = 2889.1667
$ws.Range("L80").Value = 3237.625
$ws.Range("M80").Value = -1891.1667

$ws.Range("H83").Value = 3028.55
$ws.Range("I83").Value = 2889.1667
$ws.Range("J83").Value = 3237.625
$ws.Range("K83").Value = 14445.8335
$ws.Range("L83").Value = 16188.125
$ws.Range("M83").Value = -9453.833500000001

$ws.Range("H97").Value = 569.6111
$ws.Range("I97").Value = 387.58334
$ws.Range("J97").Value = 933.6667
$ws.Range("K97").Value = 387.58334
$ws.Range("L97").Value = 933.6667
$ws.Range("M97").Value = 108.41666
$ws.Range("N97").Value = -1925.6667

$ws.Range("H132").Value = 4815.6313
$ws.Range("I132").Value = 4029.2354
$ws.Range("J132").Value = 11500
$ws.Range("K132").Value = 12087.7062
$ws.Range("L132").Value = 34500
$ws.Range("M132").Value = -9557.706200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 85789
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 85789
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 85789
$ws.Range("N12").Value = -86129

$ws.Range("H46").Value = 2619
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 2753.9167
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 2753.9167
$ws.Range("M46").Value = -812

$ws.Range("H61").Value = 5571.4
$ws.Range("I61").Value = 4439.9287
$ws.Range("J61").Value = 8211.5
$ws.Range("K61").Value = 4439.9287
$ws.Range("L61").Value = 8211.5
$ws.Range("M61").Value = -4237.9287
$ws.Range("N61").Value = -8615.5

$ws.Range("H82").Value = 1610.8823
$ws.Range("I82").Value = 1324.875
$ws.Range("J82").Value = 1865.1111
$ws.Range("K82").Value = 1324.875
$ws.Range("L82").Value = 1865.1111
$ws.Range("M82").Value = -963.875
$ws.Range("N82").Value = -2587.1111

$ws.Range("H85").Value = 1610.8823
$ws.Range("I85").Value = 1324.875
$ws.Range("J85").Value = 1865.1111
$ws.Range("K85").Value = 1324.875
$ws.Range("L85").Value = 1865.1111
$ws.Range("M85").Value = -76.875
$ws.Range("N85").Value = -4361.1111

$ws.Range("H113").Value = 5571.4
$ws.Range("I113").Value = 4439.9287
$ws.Range("J113").Value = 8211.5
$ws.Range("K113").Value = 4439.9287
$ws.Range("L113").Value = 8211.5
$ws.Range("M113").Value = -2269.9287
$ws.Range("N113").Value = -12551.5

$ws.Range("H136").Value = 8435.918
$ws.Range("I136").Value = 4396.381
$ws.Range("J136").Value = 10067.27
$ws.Range("K136").Value = 13189.143
$ws.Range("L136").Value = 30201.81
$ws.Range("M136").Value = -10639.143
$ws.Range("N136").Value = -35301.81

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 12498.333
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 12498.333
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 12498.333
$ws.Range("N41").Value = -13278.333

$ws.Range("H49").Value = 35500
$ws.Range("I49").Value = 1000
$ws.Range("J49").Value = 70000
$ws.Range("K49").Value = 1000
$ws.Range("L49").Value = 70000
$ws.Range("M49").Value = -770
$ws.Range("N49").Value = -70460

$ws.Range("H126").Value = 5489.4614
$ws.Range("I126").Value = 5969.364
$ws.Range("J126").Value = 2850
$ws.Range("K126").Value = 17908.092
$ws.Range("L126").Value = 8550
$ws.Range("M126").Value = -15438.092
$ws.Range("N126").Value = -13490

$ws.Range("H136").Value = 5666.9473
$ws.Range("I136").Value = 5181.8
$ws.Range("J136").Value = 7486.25
$ws.Range("K136").Value = 15545.4
$ws.Range("L136").Value = 22458.75
$ws.Range("M136").Value = -12995.4
$ws.Range("N136").Value = -27558.75
